$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray value in I5 (outside the table's used range)
$ws.Range("I5").ClearContents()

# Match the author's last-selected cell after the edit
$ws.Range("I5").Select()
